# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.176.87'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.26%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').Value = '3.512.86'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E3').ClearFormats()
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').Value = "'585.54"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.84%  '
$ws.Range('E5').ClearFormats()
$ws.Range('D6').Value = "'132.57"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.53%  '
$ws.Range('E6').ClearFormats()
$ws.Range('D7').Value = '3.513.58'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E7').ClearFormats()
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E8').ClearFormats()
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('E9').ClearFormats()
$ws.Range('E10').Value = '  +0.28%  '
$ws.Range('E10').ClearFormats()
$ws.Range('D11').Value = "'7.22"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('E11').ClearFormats()
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('E12').ClearFormats()
$ws.Range('D13').Value = '4.099.27'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.34%  '
$ws.Range('E13').ClearFormats()
$ws.Range('D14').Value = "'27.69"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.75%  '
$ws.Range('E14').ClearFormats()
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('E15').ClearFormats()
$ws.Range('E16').Value = '  -1.39%  '
$ws.Range('E16').ClearFormats()
$ws.Range('D17').Value = '3.518.03'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').Value = '64.244.80'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.16%  '
$ws.Range('E18').ClearFormats()
$ws.Range('E19').Value = '  +4.73%  '
$ws.Range('E19').ClearFormats()
$ws.Range('D20').Value = "'14.46"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.33%  '
$ws.Range('E20').ClearFormats()
$ws.Range('D21').Value = "'5.70"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('E21').ClearFormats()
$ws.Range('D22').Value = "'385.73"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.00%  '
$ws.Range('E22').ClearFormats()
$ws.Range('E23').Value = '  +0.36%  '
$ws.Range('E23').ClearFormats()
$ws.Range('D24').Value = '3.650.97'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.28%  '
$ws.Range('E24').ClearFormats()
$ws.Range('D25').Value = "'73.38"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.20%  '
$ws.Range('E25').ClearFormats()
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('E26').ClearFormats()
$ws.Range('D27').Value = "'0.0000115"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.69%  '
$ws.Range('E27').ClearFormats()
$ws.Range('D28').Value = "'1.57"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.88%  '
$ws.Range('E28').ClearFormats()
$ws.Range('D29').Value = "'7.54"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.36%  '
$ws.Range('E29').ClearFormats()
$ws.Range('D30').Value = "'1.00"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('E30').ClearFormats()
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('B31').ClearFormats()
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C31').ClearFormats()
$ws.Range('D31').Value = "'8.33"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.86%  '
$ws.Range('E31').ClearFormats()
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('B32').ClearFormats()
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C32').ClearFormats()
$ws.Range('D32').Value = "'2.26"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.23%  '
$ws.Range('E32').ClearFormats()
$ws.Range('D33').Value = '3.518.32'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('E33').ClearFormats()
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('E34').ClearFormats()
$ws.Range('D35').Value = "'23.82"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('E35').ClearFormats()
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('E36').ClearFormats()
$ws.Range('D37').Value = "'5.37"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.70%  '
$ws.Range('E37').ClearFormats()
$ws.Range('E38').Value = '  +0.62%  '
$ws.Range('E38').ClearFormats()
$ws.Range('D39').Value = "'6.94"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.23%  '
$ws.Range('E39').ClearFormats()
$ws.Range('D40').Value = "'161.54"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -4.88%  '
$ws.Range('E40').ClearFormats()
$ws.Range('E41').Value = '  -2.53%  '
$ws.Range('E41').ClearFormats()
$ws.Range('E42').Value = '  -0.56%  '
$ws.Range('E42').ClearFormats()
$ws.Range('D43').Value = "'26.36"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +4.35%  '
$ws.Range('E43').ClearFormats()
$ws.Range('B44').Value = 'ONDO'
$ws.Range('B44').ClearFormats()
$ws.Range('C44').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('C44').ClearFormats()
$ws.Range('D44').Value = "'1.23"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('E44').ClearFormats()
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('B45').ClearFormats()
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('C45').ClearFormats()
$ws.Range('D45').Value = "'0.998"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('E45').ClearFormats()
$ws.Range('E46').Value = '  -0.39%  '
$ws.Range('E46').ClearFormats()
$ws.Range('D47').Value = "'41.58"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.35%  '
$ws.Range('E47').ClearFormats()
$ws.Range('D48').Value = "'1.65"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.79%  '
$ws.Range('E48').ClearFormats()
$ws.Range('D49').Value = "'6.88"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('E49').ClearFormats()
$ws.Range('D50').Value = '2.438.85'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.38%  '
$ws.Range('E50').ClearFormats()
$ws.Range('B51').Value = 'VeChain'
$ws.Range('B51').ClearFormats()
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C51').ClearFormats()
$ws.Range('D51').Value = "'0.0269"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.19%  '
$ws.Range('E51').ClearFormats()
